$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("D1").Value = "description"
$ws.Range("H1").Value = "bis"

# --- Row 2 (item1) ---
$ws.Range("E2").Value = "skirmish"
$ws.Range("H2").Value = "'false"

# --- Row 3 becomes a new item2 / "Pulson grenade" row was inserted, and the
# original item2 row (Biological Transition) is pushed to row 3 with the new
# "description" / "bis" columns filled in ---
$ws.Range("A3").Value = "item2"
$ws.Range("B3").Value = "Biological Transition"
$ws.Range("C3").Value = "purple"
$ws.Range("D3").Value = "assassin,constructor,trooper"
$ws.Range("E3").Value = "skirmish"
$ws.Range("F3").Value = "xenotronics"
$ws.Range("G3").Value = "human"
$ws.Range("H3").Value = "'false"

# --- Row 4 (new) ---
$ws.Range("A4").Value = "item2"
$ws.Range("B4").Value = "Pulson grenade " + [char]8220 + "Doom D3" + [char]8221
$ws.Range("C4").Value = "purple"
$ws.Range("D4").Value = "trooper,lord commander"
$ws.Range("E4").Value = "skirmish"
$ws.Range("F4").Value = "xenotronics"
$ws.Range("G4").Value = "human"
$ws.Range("H4").Value = "'false"

$ws.Range("C5").Select() | Out-Null
